# Adds a new "2022-Q3" sheet (copied/positioned right before the existing
# "2022-Q2" sheet, with its own fund data) and updates the "总计" (summary)
# sheet with a new row for 2022-Q3, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating "2022-Q2" (so it
#    inherits identical sheet structure/formatting/outline settings),
#    placed immediately before "2022-Q2" in tab order.
# ---------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# Fill in the 2022-Q3 fund data (code, name, size, total stock position,
# position ratio, held market value, position rank).
$q3Data = @(
    @("160416", "华安标普全球石油指数（QDII-LOF）A", "2.74", "93.58", "2.47", "0.0677", 8),
    @("006282", "上投摩根欧洲动力策略股票（QDII）",   "0.41", "91.47", "3.13", "0.0128", 4),
    @("010343", "华宝英国富时100指数（QDII）A",        "0.13", "92.85", "4.23", "0.0055", 6),
    @("014982", "华安标普全球石油指数（QDII-LOF）C",   "0.22", "93.58", "2.47", "0.0054", 8),
    @("010344", "华宝英国富时100指数（QDII）C",        "0.08", "92.85", "4.23", "0.0034", 6)
)

# Columns B, D, E, F, G hold numeric-looking values that must stay TEXT
# (matching the source data's inlineStr cells), so force text format
# before writing, then clear the format residue afterwards.
$textRange = $newQ3.Range("B2:B6,D2:G6")
$textRange.NumberFormat = "@"

$r = 2
foreach ($row in $q3Data) {
    $newQ3.Cells.Item($r, 2).Value = $row[0]
    $newQ3.Cells.Item($r, 3).Value = $row[1]
    $newQ3.Cells.Item($r, 4).Value = $row[2]
    $newQ3.Cells.Item($r, 5).Value = $row[3]
    $newQ3.Cells.Item($r, 6).Value = $row[4]
    $newQ3.Cells.Item($r, 7).Value = $row[5]
    $newQ3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

$textRange.ClearFormats()

# Re-apply the bold/bordered style (matching column A's existing "s=2"
# style) to the A2:A6 index column, which ClearFormats() above did not
# touch (only B/D-G were cleared) but whose style must still read 0..4.
$idx = 0
foreach ($row in $q3Data) {
    $newQ3.Cells.Item(2 + $idx, 1).Value = $idx
    $idx = $idx + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a row for "2022-Q3" right
#    after the header, shifting all existing quarters down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$totalData = @(
    @("2022-Q3", 5, 0.09),
    @("2022-Q2", 5, 0.22),
    @("2021-Q3", 4, 0.18),
    @("2021-Q2", 3, 0.21),
    @("2021-Q1", 3, 0.26),
    @("2020-Q4", 1, 0.42)
)

$r = 2
$idx = 0
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $idx
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}

# Clean up any stray auto-formatting from the row Insert() above, then
# restore the bold/bordered "index column" style (copied from a sheet
# that was never touched) across A2:A7.
$total.Range("A2:D7").ClearFormats()
$srcQ2.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
